$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 2025Q2 row (row 28) metrics with refreshed figures
$ws.Range("C28").Value = 434
$ws.Range("D28").Value = 49
$ws.Range("E28").Value = 385
$ws.Range("F28").Value = 7.632398753894081
